$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Argoth promo card details changed (new printing appeared on market)
$ws.Range("B12").Value = "The Brothers' War Promos"
$ws.Range("C12").Value = "Normal"

# Update price column (D) values per latest market data
$ws.Range("D2").Value = 9.99
$ws.Range("D3").Value = 8.699999999999999
$ws.Range("D4").Value = 12.36
$ws.Range("D5").Value = 16.01
$ws.Range("D6").Value = 3.24
$ws.Range("D7").Value = 8.75
$ws.Range("D8").Value = 7.6
$ws.Range("D10").Value = 23.49
$ws.Range("D11").Value = 1.21
$ws.Range("D12").Value = 5.93
$ws.Range("D13").Value = 1.14
$ws.Range("D14").Value = 6.2
$ws.Range("D15").Value = 3.66
$ws.Range("D16").Value = 4.29
$ws.Range("D17").Value = 16.75
$ws.Range("D18").Value = 3.54
$ws.Range("D19").Value = 11.78
$ws.Range("D20").Value = 4.66
$ws.Range("D21").Value = 1.57
$ws.Range("D22").Value = 6.72
$ws.Range("D23").Value = 3.03
$ws.Range("D24").Value = 3.68
$ws.Range("D26").Value = 4.76
$ws.Range("D28").Value = 4.36
$ws.Range("D29").Value = 0.52
$ws.Range("D30").Value = 0.79
$ws.Range("D31").Value = 1.95
$ws.Range("D32").Value = 1.45
$ws.Range("D33").Value = 3.1
$ws.Range("D34").Value = 9.779999999999999
$ws.Range("D36").Value = 1.97
$ws.Range("D37").Value = 5.28
$ws.Range("D38").Value = 24.87
$ws.Range("D39").Value = 0.53
$ws.Range("D40").Value = 1.14
$ws.Range("D41").Value = 1.84
$ws.Range("D42").Value = 3.57
$ws.Range("D43").Value = 3.43
$ws.Range("D44").Value = 0.65
$ws.Range("D45").Value = 1.17
$ws.Range("D46").Value = 3.41
$ws.Range("D47").Value = 8.02
$ws.Range("D48").Value = 4.05
$ws.Range("D50").Value = 1.74
$ws.Range("D51").Value = 5.84
$ws.Range("D53").Value = 3.48
$ws.Range("D54").Value = 3.54
$ws.Range("D55").Value = 1.53
$ws.Range("D56").Value = 10.89
$ws.Range("D57").Value = 8.539999999999999
$ws.Range("D58").Value = 16.63
$ws.Range("D59").Value = 1.31
$ws.Range("D60").Value = 8.81
$ws.Range("D61").Value = 8.69
$ws.Range("D62").Value = 3.62
$ws.Range("D63").Value = 1.79
$ws.Range("D64").Value = 5.4
$ws.Range("D65").Value = 3.97
$ws.Range("D66").Value = 1
$ws.Range("D68").Value = 1.66

# Update active selection to reflect last edited cell
[void]$ws.Range("D39").Select()
